$wb = $excel.ActiveWorkbook

# Use an existing, already-styled sheet as the source for formatting
# (header row B1:F1 and label cell A2 both use the bold/bordered/
# centered style already present in styles.xml as style index 1).
$styleSource = $wb.Worksheets.Item(1)

# --- New sheet 1: FTNC_Demand511 -------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ftnc511 = $wb.Worksheets.Add($null, $lastSheet)
$ftnc511.Name = "FTNC_Demand511"

$ftnc511.Outline.SummaryRow = 1
$ftnc511.Outline.SummaryColumn = 1
$ftnc511.PageSetup.LeftMargin = 54
$ftnc511.PageSetup.RightMargin = 54
$ftnc511.PageSetup.TopMargin = 72
$ftnc511.PageSetup.BottomMargin = 72
$ftnc511.PageSetup.HeaderMargin = 36
$ftnc511.PageSetup.FooterMargin = 36

$ftnc511.Range("B1").Value = "In-vehicle"
$ftnc511.Range("C1").Value = "At-stop"
$ftnc511.Range("D1").Value = "Extra"
$ftnc511.Range("E1").Value = "Tardiness"
$ftnc511.Range("F1").Value = "Total"

$ftnc511.Range("A2").Value = "FTNC"
$ftnc511.Range("B2").Value = 13.30821502770083
$ftnc511.Range("C2").Value = 182.7941640821028
$ftnc511.Range("D2").Value = 0
$ftnc511.Range("E2").Value = 53.38723945807137
$ftnc511.Range("F2").Value = 249.4896185678749

$styleSource.Range("B1:F1").Copy()
$ftnc511.Range("B1:F1").PasteSpecial(-4122)
$styleSource.Range("A2").Copy()
$ftnc511.Range("A2").PasteSpecial(-4122)

# --- New sheet 2: FTHC_Demand5 ----------------------------------------
$fthc5 = $wb.Worksheets.Add($null, $ftnc511)
$fthc5.Name = "FTHC_Demand5"

$fthc5.Outline.SummaryRow = 1
$fthc5.Outline.SummaryColumn = 1
$fthc5.PageSetup.LeftMargin = 54
$fthc5.PageSetup.RightMargin = 54
$fthc5.PageSetup.TopMargin = 72
$fthc5.PageSetup.BottomMargin = 72
$fthc5.PageSetup.HeaderMargin = 36
$fthc5.PageSetup.FooterMargin = 36

$fthc5.Range("B1").Value = "In-vehicle"
$fthc5.Range("C1").Value = "At-stop"
$fthc5.Range("D1").Value = "Extra"
$fthc5.Range("E1").Value = "Tardiness"
$fthc5.Range("F1").Value = "Total"

$fthc5.Range("A2").Value = "FTHC"
$fthc5.Range("B2").Value = 12.44083242684707
$fthc5.Range("C2").Value = 187.3018996736528
$fthc5.Range("D2").Value = 0
$fthc5.Range("E2").Value = 8.239537621283734
$fthc5.Range("F2").Value = 207.9822697217838

$styleSource.Range("B1:F1").Copy()
$fthc5.Range("B1:F1").PasteSpecial(-4122)
$styleSource.Range("A2").Copy()
$fthc5.Range("A2").PasteSpecial(-4122)

$fthc5.Range("A1").Select()
